# Natmi following Dr Hou advice
# The LR-pairs sheet was recomputed with updated ligand/receptor-expressing
# cell counts (E, K: 1 -> 3) and refreshed average-expression inputs
# (G per sending-cluster row-block, M per target-cluster column). Every
# other touched column (H, I, J, N, O, P, Q, R, S, T) is a pure function of
# those two counts + two expression inputs, so we recompute them instead of
# hardcoding the ~150 changed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

# New ligand-expressing / receptor-expressing cell counts for every data row.
$newCount = 3

# New "average expression" inputs, keyed by the existing (old) average so we
# don't have to hardcode row numbers -- every row in a sending-cluster block
# shares one G value, every row in a target-cluster column shares one M
# value, and both stay uniquely identifiable by their old value.
$newG = @{
    "12.0339014791582" = 14.05492
    "5.45867822632938"  = 5.498465333333333
    "3.78791222965858"  = 4.576766333333333
    "1.74312255952082"  = 2.239594666666667
}
$newM = @{
    "1.49714132510024" = 1.499230333333333
    "3.29685486510907" = 3.340648
    "3.26331753920787" = 4.215231
    "3.57237684716374" = 3.615388333333333
}

# Pass 1: write the counts (E, K) and the refreshed average-expression
# inputs (G, M); derive the matching totals (H, N) while we're at it.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("E$r").Value = $newCount
    $ws.Range("K$r").Value = $newCount

    $oldG = $ws.Range("G$r").Value2
    $oldM = $ws.Range("M$r").Value2

    $gVal = $newG[[string]$oldG]
    $mVal = $newM[[string]$oldM]

    $ws.Range("G$r").Value = $gVal
    $ws.Range("M$r").Value = $mVal

    $ws.Range("H$r").Value = $gVal * $newCount
    $ws.Range("N$r").Value = $mVal * $newCount
}

# Pass 2: ligand/receptor derived specificity (I/J, O/P) -- each row's
# average expression divided by the sum of averages across the sheet's
# sending clusters (for G) / target clusters (for M). Re-read back through
# Value2 so the sums reflect what's now on the sheet.
$sumG = 0.0
$sumM = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $sumG += $ws.Range("G$r").Value2
    $sumM += $ws.Range("M$r").Value2
}
# Each distinct G repeats 4x (one per target cluster) and each distinct M
# repeats 4x (one per sending cluster), so divide the running totals back
# down to a "per distinct value" sum.
$sumG = $sumG / 4.0
$sumM = $sumM / 4.0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $gVal = $ws.Range("G$r").Value2
    $mVal = $ws.Range("M$r").Value2

    $iVal = $gVal / $sumG
    $oVal = $mVal / $sumM

    $ws.Range("I$r").Value = $iVal
    $ws.Range("J$r").Value = $iVal
    $ws.Range("O$r").Value = $oVal
    $ws.Range("P$r").Value = $oVal
}

# Pass 3: edge weights (Q = G*M average weight, R = H*N total weight) and
# their sheet-wide derived specificities (S/T).
$sumQ = 0.0
$sumR = 0.0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $qVal = $ws.Range("G$r").Value2 * $ws.Range("M$r").Value2
    $rVal = $ws.Range("H$r").Value2 * $ws.Range("N$r").Value2

    $ws.Range("Q$r").Value = $qVal
    $ws.Range("R$r").Value = $rVal

    $sumQ += $qVal
    $sumR += $rVal
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $qVal = $ws.Range("Q$r").Value2
    $rVal = $ws.Range("R$r").Value2

    $ws.Range("S$r").Value = $qVal / $sumQ
    $ws.Range("T$r").Value = $rVal / $sumR
}
